# Insert a new data row at row 65 (shifts existing rows 65-152 down to 66-153),
# then populate the new row 65 with the same record as the original row 65 had,
# except for a new Fecha (date) and Volumen value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 65; this shifts rows 65..152 down to 66..153
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the original row's data, except D (Fecha) and J (Volumen)
$ws.Range("A65").Value = 1
$ws.Range("B65").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C65").Value = 'Arica y Parinacota'
$ws.Range("D65").Value = 44413
$ws.Range("E65").Value = 15
$ws.Range("F65").Value = 100114013
$ws.Range("G65").Value = 'Zanahoria'
$ws.Range("H65").Value = 'Sin especificar'
$ws.Range("I65").Value = 'Primera'
$ws.Range("J65").Value = 70
$ws.Range("K65").Value = 11000
$ws.Range("L65").Value = 12000
$ws.Range("M65").Value = 11500
$ws.Range("N65").Value = '$/saco 25 kilos'
$ws.Range("O65").Value = 'Valle de Camiña'
$ws.Range("P65").Value = 460
$ws.Range("Q65").Value = 25
$ws.Range("R65").Value = 'Hortaliza'

# Make sure the date cell keeps the expected date number format
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
